# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The 16-row "Periodo Mora" / "Valor Mora" block (rows 16-31) is re-sorted
# from descending (2002 .. 1810) to ascending (1810 .. 2002) order; since
# every row shares the same worker (B/C/D) and the same base salary (G),
# the net effect is that the Periodo Mora label and Valor Mora amount for
# each row position are replaced with the values from the mirrored
# (16+31-row) position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("1810","1811","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002")
$valores  = @(31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,17708)

for ($i = 0; $i -lt 16; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 5).Value = $periodos[$i]
    $ws.Cells.Item($r, 6).Value = $valores[$i]
}
